# The "Recorded By" column (G) lists the two recorders separated by a comma.
# The author order needs to be swapped from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System" everywhere it appears on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G$r")
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
